$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(538500, "Glyph not found"),
    @(538501, "Non Message Handelr for the message"),
    @(538502, "Unknown item command:  for custome menu"),
    @(538503, "Unknown Menu Option for Application Menu"),
    @(538504, "Unsupported Config Version"),
    @(538505, $null),
    @(538506, $null),
    @(538507, $null),
    @(538508, $null),
    @(538509, $null),
    @(538510, $null),
    @(538511, $null),
    @(538512, $null),
    @(538513, $null),
    @(538514, $null),
    @(538515, $null),
    @(538516, $null),
    @(538517, $null),
    @(538518, $null),
    @(538519, $null),
    @(538520, $null),
    @(538521, $null)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    if ($item[1] -ne $null) {
        $ws.Cells.Item($row, 2).Value = $item[1]
    }
    $row = $row + 1
}

$ws.Columns.Item(2).ColumnWidth = 41.7109375

$ws.Range("B10").Select() | Out-Null
